# Fruta / hortaliza, semanal
# Rewrite rows 4-10 (columns D, K, L, M, N, O, P, Q, R, S, T) with the
# updated weekly data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44545
$ws.Range("K4").Value = "Castle Brite"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 19000
$ws.Range("P4").Value = 18500
$ws.Range("Q4").Value = "`$/caja 15 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1233
$ws.Range("T4").Value = 15

# Row 5
$ws.Range("D5").Value = 44545
$ws.Range("K5").Value = "Castle Brite"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 17000
$ws.Range("P5").Value = 17000
$ws.Range("Q5").Value = "`$/caja 15 kilos"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1133
$ws.Range("T5").Value = 15

# Row 6
$ws.Range("D6").Value = 44189
$ws.Range("K6").Value = "Dina"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = "`$/caja 15 kilos granel"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1033
$ws.Range("T6").Value = 15

# Row 7
$ws.Range("D7").Value = 44189
$ws.Range("K7").Value = "Dina"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 14000
$ws.Range("Q7").Value = "`$/caja 15 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 933
$ws.Range("T7").Value = 15

# Row 8 is unchanged.

# Row 9
$ws.Range("D9").Value = 44159
$ws.Range("K9").Value = "Castle Brite"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("Q9").Value = "`$/caja 15 kilos"
$ws.Range("R9").Value = "Región Metropolitana"
$ws.Range("S9").Value = 967
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("D10").Value = 44579
$ws.Range("K10").Value = "Modesto"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 180
$ws.Range("N10").Value = 13000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 13444
$ws.Range("Q10").Value = "`$/caja 18 kilos"
$ws.Range("R10").Value = "Región Metropolitana"
$ws.Range("S10").Value = 747
$ws.Range("T10").Value = 18
